$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("L2").Value = "[58.88844173295023, 68.05615866710345]"
$ws.Range("T2").Value = "[46.3963224571107, 52.68256152407314]"
$ws.Range("L3").Value = "[59.74163852128687, 67.468116364321]"
$ws.Range("T3").Value = "[48.20381468286139, 52.30516770209411]"
